$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.680.26"
$ws.Range("E2").Value = "  -6.13%  "
$ws.Range("D3").Value = "2.610.29"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "302.17"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("D6").Value = "96.68"
$ws.Range("E6").Value = "  -3.79%  "
$ws.Range("D7").Value = "0.579"
$ws.Range("E7").Value = "  -3.97%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "0.561"
$ws.Range("E9").Value = "  -3.17%  "
$ws.Range("D10").Value = "36.97"
$ws.Range("E10").Value = "  -6.29%  "
$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  -3.48%  "
$ws.Range("D12").Value = "7.84"
$ws.Range("E12").Value = "  -4.07%  "
$ws.Range("D13").Value = "3.002.48"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "2.602.39"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "0.894"
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("D17").Value = "14.41"
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("D18").Value = "43.626.81"
$ws.Range("E18").Value = "  -6.33%  "
$ws.Range("D19").Value = "6.71"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("E20").Value = "  -3.71%  "
$ws.Range("D21").Value = "12.37"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").Value = "73.22"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").Value = "266.44"
$ws.Range("E23").Value = "  -4.27%  "
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("D25").Value = "2.22"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").Value = "29.38"
$ws.Range("E26").Value = "  +3.16%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "10.32"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "37.76"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.16"
$ws.Range("E30").Value = "  -6.54%  "
$ws.Range("D31").Value = "6.09"
$ws.Range("E31").Value = "  -5.42%  "
$ws.Range("D32").Value = "3.64"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "152.39"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "0.0815"
$ws.Range("E37").Value = "  -4.30%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "0.121"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "24.31"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").Value = "16.76"
$ws.Range("E40").Value = "  +3.54%  "
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").Value = "0.0315"
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("D43").Value = "3.88"
$ws.Range("E43").Value = "  -5.27%  "
$ws.Range("D44").Value = "2.044.58"
$ws.Range("E44").Value = "  -4.11%  "
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "88.50"
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("D47").Value = "9.14"
$ws.Range("E47").Value = "  -4.30%  "
$ws.Range("D48").Value = "1.63"
$ws.Range("E48").Value = "  +4.25%  "
$ws.Range("D49").Value = "2.855.90"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "106.36"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").Value = "0.192"
$ws.Range("E51").Value = "  -4.36%  "
